$d = $word.ActiveDocument

# Locate the paragraph that ends with the final bibliography line and the
# paragraph holding the closing "(c) 2020 ..." footer line, then delete the
# empty paragraph plus the two footer paragraphs that follow the
# bibliography entry ("Ver no Jupiter ..." and "(c) 2020 ...").
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*23.ed. S*o Paulo: Cortez, 2009.*") {
        $anchor = $p
    }
}

$footerEnd = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Contact: luizeleno@usp.br*") {
        $footerEnd = $p
    }
}

$r = $d.Range($anchor.Range.End, $footerEnd.Range.End)
$r.Delete()
